$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.104402212121897
$ws.Range("C2").Value = 0.02741930477330357
$ws.Range("D2").Value = 0.4930542803894156
$ws.Range("E2").Value = 0.145548477137325
$ws.Range("G2").Value = 3.045865065099576
$ws.Range("H2").Value = 2.366088401623443
$ws.Range("I2").Value = 2.435079019361133
$ws.Range("J2").Value = 0.04433847406113811
$ws.Range("K2").Value = 0.8210316687575414
$ws.Range("L2").Value = 0.5301789097841692
$ws.Range("B3").Value = 1.080202246912137
$ws.Range("C3").Value = 0.02385532807493007
$ws.Range("D3").Value = 0.4908237796530841
$ws.Range("E3").Value = 0.1456024328986985
$ws.Range("G3").Value = 3.031081629606973
$ws.Range("H3").Value = 2.364444310796955
$ws.Range("I3").Value = 2.430583480470979
$ws.Range("J3").Value = 0.04442363888641854
$ws.Range("K3").Value = 0.7945556817903423
$ws.Range("L3").Value = 0.5267605603038703
$ws.Range("B4").Value = 1.066005558229534
$ws.Range("C4").Value = 0.02166690258246717
$ws.Range("D4").Value = 0.4896567111469778
$ws.Range("E4").Value = 0.1456821694627113
$ws.Range("G4").Value = 3.023078369973078
$ws.Range("H4").Value = 2.364135242036809
$ws.Range("I4").Value = 2.428562245600055
$ws.Range("J4").Value = 0.04447939325268324
$ws.Range("K4").Value = 0.778823655251017
$ws.Range("L4").Value = 0.5248987652203425
$ws.Range("B5").Value = 1.060387193872884
$ws.Range("C5").Value = 0.02077503125183
$ws.Range("D5").Value = 0.4892321078726667
$ws.Range("E5").Value = 0.1457264051464087
$ws.Range("G5").Value = 3.020086846461183
$ws.Range("H5").Value = 2.364185480510287
$ws.Range("I5").Value = 2.42792439153623
$ws.Range("J5").Value = 0.04450298653988094
$ws.Range("K5").Value = 0.7725448126576566
$ws.Range("L5").Value = 0.5241997668671416
$ws.Range("B6").Value = 1.059464360780453
$ws.Range("C6").Value = 0.02062693136797122
$ws.Range("D6").Value = 0.4891646844992437
$ws.Range("E6").Value = 0.1457344603257305
$ws.Range("G6").Value = 3.019606401330634
$ws.Range("H6").Value = 2.364204465697355
$ws.Range("I6").Value = 2.427829698533124
$ws.Range("J6").Value = 0.04450695697412543
$ws.Range("K6").Value = 0.7715102004174526
$ws.Range("L6").Value = 0.5240873067711931
$ws.Range("B7").Value = 1.065929110706662
$ws.Range("C7").Value = 0.02165487482770345
$ws.Range("D7").Value = 0.4896507782560917
$ws.Range("E7").Value = 0.1456827184659506
$ws.Range("G7").Value = 3.023036932760988
$ws.Range("H7").Value = 2.364135206125013
$ws.Range("I7").Value = 2.428552890956567
$ws.Range("J7").Value = 0.0444797079025645
$ws.Range("K7").Value = 0.7787384414552321
$ws.Range("L7").Value = 0.5248890964590487
$ws.Range("B8").Value = 1.095920847787141
$ws.Range("C8").Value = 0.02619044870378673
$ws.Range("D8").Value = 0.4922432241422001
$ws.Range("E8").Value = 0.145557420938248
$ws.Range("G8").Value = 3.040544760526728
$ws.Range("H8").Value = 2.365376182066171
$ws.Range("I8").Value = 2.433375549459974
$ws.Range("J8").Value = 0.04436712149286892
$ws.Range("K8").Value = 0.8117940523562481
$ws.Range("L8").Value = 0.5289511019017539
$ws.Range("B9").Value = 1.159976480278402
$ws.Range("C9").Value = 0.03508579988839244
$ws.Range("D9").Value = 0.4989310276934447
$ws.Range("E9").Value = 0.1456806425502108
$ws.Range("G9").Value = 3.083409325445047
$ws.Range("H9").Value = 2.37336642241641
$ws.Range("I9").Value = 2.448700238744038
$ws.Range("J9").Value = 0.04417371819064941
$ws.Range("K9").Value = 0.880770773128063
$ws.Range("L9").Value = 0.5387953321445309
$ws.Range("B10").Value = 1.210224512754849
$ws.Range("C10").Value = 0.0416254625330339
$ws.Range("D10").Value = 0.5048202843528458
$ws.Range("E10").Value = 0.1459950638270335
$ws.Range("G10").Value = 3.120125408934058
$ws.Range("H10").Value = 2.382626349754048
$ws.Range("I10").Value = 2.463544619096837
$ws.Range("J10").Value = 0.04404818112996356
$ws.Range("K10").Value = 0.9339804275618064
$ws.Range("L10").Value = 0.5471713426016009
$ws.Range("B11").Value = 1.233774128553279
$ws.Range("C11").Value = 0.04460225718433719
$ws.Range("D11").Value = 0.5077109719046859
$ws.Range("E11").Value = 0.1461865160173161
$ws.Range("G11").Value = 3.137967984384517
$ws.Range("H11").Value = 2.387575592850169
$ws.Range("I11").Value = 2.471078289548217
$ws.Range("J11").Value = 0.04399463779657697
$ws.Range("K11").Value = 0.9587372395302793
$ws.Range("L11").Value = 0.5512298354336878
$ws.Range("B12").Value = 1.242790941841292
$ws.Range("C12").Value = 0.04572981270374044
$ws.Range("D12").Value = 0.5088359815929522
$ws.Range("E12").Value = 0.1462659575130871
$ws.Range("G12").Value = 3.144888769609679
$ws.Range("H12").Value = 2.389555704920042
$ws.Range("I12").Value = 2.474043487940833
$ws.Range("J12").Value = 0.04397487266592259
$ws.Range("K12").Value = 0.9681911984489204
$ws.Range("L12").Value = 0.5528023221618525
$ws.Range("B13").Value = 1.240844608556387
$ws.Range("C13").Value = 0.04548695912225753
$ws.Range("D13").Value = 0.5085923409631761
$ws.Range("E13").Value = 0.1462485399045939
$ws.Range("G13").Value = 3.143390948743871
$ws.Range("H13").Value = 2.389124540856699
$ws.Range("I13").Value = 2.47339988099489
$ws.Range("J13").Value = 0.043979106763683
$ws.Range("K13").Value = 0.9661516049048373
$ws.Range("L13").Value = 0.5524620757073251
$ws.Range("B14").Value = 1.234513963422813
$ws.Range("C14").Value = 0.04469501545813159
$ws.Range("D14").Value = 0.5078029188506861
$ws.Range("E14").Value = 0.1461929126891555
$ws.Range("G14").Value = 3.138534069668907
$ws.Range("H14").Value = 2.387736374740086
$ws.Range("I14").Value = 2.471319985917845
$ws.Range("J14").Value = 0.04399300148606811
$ws.Range("K14").Value = 0.9595134384893527
$ws.Range("L14").Value = 0.5513584912975347
$ws.Range("B15").Value = 1.23064915377617
$ws.Range("C15").Value = 0.04420996837562541
$ws.Range("D15").Value = 0.5073233280807443
$ws.Range("E15").Value = 0.1461597429661623
$ws.Range("G15").Value = 3.135580479933083
$ws.Range("H15").Value = 2.38689987907415
$ws.Range("I15").Value = 2.470060624940871
$ws.Range("J15").Value = 0.04400157883834499
$ws.Range("K15").Value = 0.9554576638101366
$ws.Range("L15").Value = 0.5506871512323812
$ws.Range("B16").Value = 1.208699374760101
$ws.Range("C16").Value = 0.04143096295705107
$ws.Range("D16").Value = 0.5046356242862089
$ws.Range("E16").Value = 0.145983524596371
$ws.Range("G16").Value = 3.118982318346156
$ws.Range("H16").Value = 2.382317730929429
$ws.Range("I16").Value = 2.463067995436248
$ws.Range("J16").Value = 0.04405175186515109
$ws.Range("K16").Value = 0.9323735906791626
$ws.Range("L16").Value = 0.546911099215194
$ws.Range("B17").Value = 1.195410758805053
$ws.Range("C17").Value = 0.03972663968345103
$ws.Range("D17").Value = 0.5030409641924081
$ws.Range("E17").Value = 0.145887807599383
$ws.Range("G17").Value = 3.109092086102152
$ws.Range("H17").Value = 2.379695444595797
$ws.Range("I17").Value = 2.458978297251491
$ws.Range("J17").Value = 0.04408344285384747
$ws.Range("K17").Value = 0.9183533540958422
$ws.Range("L17").Value = 0.5446581358536946
$ws.Range("B18").Value = 1.187832609746977
$ws.Range("C18").Value = 0.03874653159920172
$ws.Range("D18").Value = 0.5021436817487199
$ws.Range("E18").Value = 0.1458373124579744
$ws.Range("G18").Value = 3.103510797831404
$ws.Range("H18").Value = 2.378256532104956
$ws.Range("I18").Value = 2.456699507519218
$ws.Range("J18").Value = 0.04410200625845384
$ws.Range("K18").Value = 0.9103412042927914
$ws.Range("L18").Value = 0.5433856569942037
$ws.Range("B19").Value = 1.185277972604155
$ws.Range("C19").Value = 0.03841471246776962
$ws.Range("D19").Value = 0.5018433012326824
$ws.Range("E19").Value = 0.1458209992543971
$ws.Range("G19").Value = 3.101639493660002
$ws.Range("H19").Value = 2.377781253823059
$ws.Range("I19").Value = 2.455940569459514
$ws.Range("J19").Value = 0.04410834920499429
$ws.Range("K19").Value = 0.907637351735076
$ws.Range("L19").Value = 0.5429588327499886
$ws.Range("B20").Value = 1.196818618480762
$ws.Range("C20").Value = 0.03990804955810745
$ws.Range("D20").Value = 0.5032086569035243
$ws.Range("E20").Value = 0.1458975251706711
$ws.Range("G20").Value = 3.110133810227325
$ws.Range("H20").Value = 2.379967413361612
$ws.Range("I20").Value = 2.459406045480137
$ws.Range("J20").Value = 0.04408003457627707
$ws.Range("K20").Value = 0.9198404602821881
$ws.Range("L20").Value = 0.5448955498446537
$ws.Range("B21").Value = 1.236370741032971
$ws.Range("C21").Value = 0.04492761983613036
$ws.Range("D21").Value = 0.5080339675573669
$ws.Range("E21").Value = 0.1462090634868609
$ws.Range("G21").Value = 3.139956194948098
$ws.Range("H21").Value = 2.388141237697994
$ws.Range("I21").Value = 2.471927851509378
$ws.Range("J21").Value = 0.04398890642791509
$ws.Range("K21").Value = 0.9614610839991826
$ws.Range("L21").Value = 0.5516816745107036
$ws.Range("B22").Value = 1.262797720913454
$ws.Range("C22").Value = 0.04821001282461168
$ws.Range("D22").Value = 0.5113645721055207
$ws.Range("E22").Value = 0.1464531300136827
$ws.Range("G22").Value = 3.160403972569668
$ws.Range("H22").Value = 2.394100803492108
$ws.Range("I22").Value = 2.480766534563443
$ws.Range("J22").Value = 0.04393232400314018
$ws.Range("K22").Value = 0.9891234801913527
$ws.Range("L22").Value = 0.5563244101956002
$ws.Range("B23").Value = 1.248640430238993
$ws.Range("C23").Value = 0.0464579589187224
$ws.Range("D23").Value = 0.5095707924751878
$ws.Range("E23").Value = 0.1463191717679919
$ws.Range("G23").Value = 3.149402957434518
$ws.Range("H23").Value = 2.390863575257697
$ws.Range("I23").Value = 2.475989210704157
$ws.Range("J23").Value = 0.04396225152652278
$ws.Range("K23").Value = 0.9743174330262718
$ws.Range("L23").Value = 0.5538275204300618
$ws.Range("B24").Value = 1.196181933048678
$ws.Range("C24").Value = 0.03982603499864013
$ws.Range("D24").Value = 0.5031327822577651
$ws.Range("E24").Value = 0.1458931177329852
$ws.Range("G24").Value = 3.109662520499825
$ws.Range("H24").Value = 2.379844242343154
$ws.Range("I24").Value = 2.459212434873606
$ws.Range("J24").Value = 0.04408157438741211
$ws.Range("K24").Value = 0.9191679892005595
$ws.Range("L24").Value = 0.5447881440185398
$ws.Range("B25").Value = 1.142087531055438
$ws.Range("C25").Value = 0.032678885283147
$ws.Range("D25").Value = 0.4969503013134045
$ws.Range("E25").Value = 0.145607917472848
$ws.Range("G25").Value = 3.070897644305404
$ws.Range("H25").Value = 2.37060977030211
$ws.Range("I25").Value = 2.443925492004752
$ws.Range("J25").Value = 0.04422312181134869
$ws.Range("K25").Value = 0.8616660085305909
$ws.Range("L25").Value = 0.5359312778147114
